$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest scrape: cell -> new text value (coinranking.com snapshot refresh).
$updates = @(
    @{ Cell = "D2"; Value = '27.968.18' }
    @{ Cell = "E2"; Value = '  +0.81%  ' }
    @{ Cell = "D3"; Value = '1.768.76' }
    @{ Cell = "E3"; Value = '  -0.32%  ' }
    @{ Cell = "D4"; Value = '1.001' }
    @{ Cell = "E4"; Value = '  -0.08%  ' }
    @{ Cell = "D5"; Value = '328.70' }
    @{ Cell = "E5"; Value = '  +0.42%  ' }
    @{ Cell = "E6"; Value = '  -0.08%  ' }
    @{ Cell = "D7"; Value = '0.4657' }
    @{ Cell = "E7"; Value = '  +1.63%  ' }
    @{ Cell = "D8"; Value = '0.3517' }
    @{ Cell = "E8"; Value = '  -1.88%  ' }
    @{ Cell = "D9"; Value = '43.48' }
    @{ Cell = "E9"; Value = '  +3.87%  ' }
    @{ Cell = "D10"; Value = '0.07392' }
    @{ Cell = "E10"; Value = '  -1.31%  ' }
    @{ Cell = "D11"; Value = '1.084' }
    @{ Cell = "E11"; Value = '  -1.91%  ' }
    @{ Cell = "D12"; Value = '1.000' }
    @{ Cell = "E12"; Value = '  -0.08%  ' }
    @{ Cell = "D13"; Value = '20.62' }
    @{ Cell = "E13"; Value = '  -1.10%  ' }
    @{ Cell = "D14"; Value = '6.014' }
    @{ Cell = "E14"; Value = '  -0.54%  ' }
    @{ Cell = "D15"; Value = '7.200' }
    @{ Cell = "E15"; Value = '  -0.35%  ' }
    @{ Cell = "D16"; Value = '1.766.91' }
    @{ Cell = "E16"; Value = '  -0.65%  ' }
    @{ Cell = "D17"; Value = '92.25' }
    @{ Cell = "E17"; Value = '  -1.59%  ' }
    @{ Cell = "D18"; Value = '0.00001056' }
    @{ Cell = "E18"; Value = '  -0.41%  ' }
    @{ Cell = "D19"; Value = '0.06426' }
    @{ Cell = "E19"; Value = '  -0.12%  ' }
    @{ Cell = "E20"; Value = '  -0.07%  ' }
    @{ Cell = "D21"; Value = '16.95' }
    @{ Cell = "E21"; Value = '  -0.83%  ' }
    @{ Cell = "D22"; Value = '5.796' }
    @{ Cell = "E22"; Value = '  -0.26%  ' }
    @{ Cell = "D23"; Value = '28.003.68' }
    @{ Cell = "E23"; Value = '  +0.81%  ' }
    @{ Cell = "D24"; Value = '11.15' }
    @{ Cell = "E24"; Value = '  -1.55%  ' }
    @{ Cell = "D25"; Value = '2.158' }
    @{ Cell = "E25"; Value = '  +3.57%  ' }
    @{ Cell = "D26"; Value = '163.98' }
    @{ Cell = "E26"; Value = '  -0.35%  ' }
    @{ Cell = "D27"; Value = '20.03' }
    @{ Cell = "E27"; Value = '  -1.22%  ' }
    @{ Cell = "D28"; Value = '1.970.15' }
    @{ Cell = "E28"; Value = '  -0.47%  ' }
    @{ Cell = "D29"; Value = '2.196' }
    @{ Cell = "E29"; Value = '  +1.21%  ' }
    @{ Cell = "D30"; Value = '123.34' }
    @{ Cell = "E30"; Value = '  -1.87%  ' }
    @{ Cell = "D31"; Value = '1.077' }
    @{ Cell = "E31"; Value = '  -2.37%  ' }
    @{ Cell = "D32"; Value = '0.09325' }
    @{ Cell = "E32"; Value = '  +1.17%  ' }
    @{ Cell = "D33"; Value = '3.655' }
    @{ Cell = "E33"; Value = '  -0.41%  ' }
    @{ Cell = "D34"; Value = '5.555' }
    @{ Cell = "E34"; Value = '  +0.19%  ' }
    @{ Cell = "E35"; Value = '  -1.56%  ' }
    @{ Cell = "B36"; Value = 'VeChain' }
    @{ Cell = "C36"; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = "D36"; Value = '0.02271' }
    @{ Cell = "E36"; Value = '  -1.15%  ' }
    @{ Cell = "B37"; Value = 'Hedera' }
    @{ Cell = "C37"; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = "D37"; Value = '0.06105' }
    @{ Cell = "E37"; Value = '  -1.28%  ' }
    @{ Cell = "D38"; Value = '0.2067' }
    @{ Cell = "E38"; Value = '  -1.10%  ' }
    @{ Cell = "D39"; Value = '4.923' }
    @{ Cell = "E39"; Value = '  -0.81%  ' }
    @{ Cell = "D40"; Value = '1.448' }
    @{ Cell = "E40"; Value = '  +4.22%  ' }
    @{ Cell = "B41"; Value = 'TrustWalletToken' }
    @{ Cell = "C41"; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = "D41"; Value = '1.191' }
    @{ Cell = "E41"; Value = '  +0.10%  ' }
    @{ Cell = "B42"; Value = 'TheSandbox' }
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' }
    @{ Cell = "D42"; Value = '0.6157' }
    @{ Cell = "E42"; Value = '  -2.62%  ' }
    @{ Cell = "D43"; Value = '7.764' }
    @{ Cell = "E43"; Value = '  -0.46%  ' }
    @{ Cell = "E44"; Value = '  -1.47%  ' }
    @{ Cell = "D45"; Value = '3.744' }
    @{ Cell = "E45"; Value = '  -0.08%  ' }
    @{ Cell = "D46"; Value = '0.5799' }
    @{ Cell = "E46"; Value = '  -1.86%  ' }
    @{ Cell = "D47"; Value = '123.97' }
    @{ Cell = "E47"; Value = '  +1.04%  ' }
    @{ Cell = "D48"; Value = '1.934' }
    @{ Cell = "E48"; Value = '  -0.94%  ' }
    @{ Cell = "B49"; Value = 'EOS' }
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos' }
    @{ Cell = "D49"; Value = '1.126' }
    @{ Cell = "E49"; Value = '  -1.21%  ' }
    @{ Cell = "B50"; Value = 'Cronos' }
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = "D50"; Value = '0.06815' }
    @{ Cell = "E50"; Value = '  -1.58%  ' }
    @{ Cell = "D51"; Value = '72.10' }
    @{ Cell = "E51"; Value = '  -0.45%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Price/Coin/Link/Volume columns are plain text in this sheet (e.g. "1.001",
    # "27.938.71" use dots as thousand separators) - force text format first so
    # Excel does not reinterpret numeric-looking strings as numbers/dates, then
    # drop back to the default style so no formatting change is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
